# Weekly update: insert two new daily price records (Primera / Segunda quality)
# for "Acelga" at "Vega Monumental Concepción" dated 2023-10-19 (serial 45218),
# pushing the existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 421, shifting all the
# existing data (rows 421..477) down to rows 423..479.
$ws.Rows.Item(421).Insert()
$ws.Rows.Item(421).Insert()

# New row 421: Primera quality
$ws.Cells.Item(421, 1).Value = 11
$ws.Cells.Item(421, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(421, 3).Value = "Bíobío"
$ws.Cells.Item(421, 4).Value = 45218
$ws.Cells.Item(421, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(421, 5).Value = 8
$ws.Cells.Item(421, 6).Value = 100112009
$ws.Cells.Item(421, 7).Value = "Acelga"
$ws.Cells.Item(421, 8).Value = "Sin especificar"
$ws.Cells.Item(421, 9).Value = "Primera"
$ws.Cells.Item(421, 10).Value = 200
$ws.Cells.Item(421, 11).Value = 600
$ws.Cells.Item(421, 12).Value = 700
$ws.Cells.Item(421, 13).Value = 650
$ws.Cells.Item(421, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(421, 15).Value = "Región de Ñuble"
$ws.Cells.Item(421, 16).Value = 650
$ws.Cells.Item(421, 17).Value = 1
$ws.Cells.Item(421, 18).Value = "Hortaliza"

# New row 422: Segunda quality
$ws.Cells.Item(422, 1).Value = 11
$ws.Cells.Item(422, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(422, 3).Value = "Bíobío"
$ws.Cells.Item(422, 4).Value = 45218
$ws.Cells.Item(422, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(422, 5).Value = 8
$ws.Cells.Item(422, 6).Value = 100112009
$ws.Cells.Item(422, 7).Value = "Acelga"
$ws.Cells.Item(422, 8).Value = "Sin especificar"
$ws.Cells.Item(422, 9).Value = "Segunda"
$ws.Cells.Item(422, 10).Value = 100
$ws.Cells.Item(422, 11).Value = 500
$ws.Cells.Item(422, 12).Value = 500
$ws.Cells.Item(422, 13).Value = 500
$ws.Cells.Item(422, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(422, 15).Value = "Región de Ñuble"
$ws.Cells.Item(422, 16).Value = 500
$ws.Cells.Item(422, 17).Value = 1
$ws.Cells.Item(422, 18).Value = "Hortaliza"

Write-Host "Inserted two new rows for Acelga (Primera/Segunda) at row 421-422"
